$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 542.5
$ws.Range("I19").Value = 607.2308
$ws.Range("J19").Value = 374.2
$ws.Range("K19").Value = 607.2308
$ws.Range("L19").Value = 374.2
$ws.Range("M19").Value = -432.2308
$ws.Range("N19").Value = -724.2
$ws.Range("H38").Value = 1668.2858
$ws.Range("I38").Value = 153.8
$ws.Range("K38").Value = 461.4
$ws.Range("M38").Value = -89.40000000000003
$ws.Range("H52").Value = 462.5
$ws.Range("I52").Value = 462.5
$ws.Range("K52").Value = 1387.5
$ws.Range("M52").Value = -1227.5
$ws.Range("H113").Value = 1900.3334
$ws.Range("I113").Value = 2296.6667
$ws.Range("J113").Value = 1504
$ws.Range("K113").Value = 2296.6667
$ws.Range("L113").Value = 1504
$ws.Range("M113").Value = 957.3332999999998
$ws.Range("N113").Value = -8012
$ws.Range("H135").Value = 62502124
$ws.Range("I135").Value = 66668800
$ws.Range("K135").Value = 600019200
$ws.Range("M135").Value = -600016665
$ws.Range("H137").Value = 1117829.4
$ws.Range("I137").Value = 2010
$ws.Range("J137").Value = 1505940.5
$ws.Range("K137").Value = 6030
$ws.Range("L137").Value = 4517821.5
$ws.Range("M137").Value = -3480
$ws.Range("N137").Value = -4522921.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1139.5
$ws.Range("I2").Value = 1082.7028
$ws.Range("K2").Value = 1082.7028
$ws.Range("M2").Value = -969.7028
$ws.Range("H32").Value = 6416142.5
$ws.Range("I32").Value = 6761082.5
$ws.Range("K32").Value = 6761082.5
$ws.Range("M32").Value = -6760795.5
$ws.Range("H61").Value = 21657.691
$ws.Range("I61").Value = 25222.143
$ws.Range("K61").Value = 25222.143
$ws.Range("M61").Value = -25010.143
$ws.Range("H74").Value = 2780.831
$ws.Range("I74").Value = 2503
$ws.Range("K74").Value = 2503
$ws.Range("M74").Value = -1629
$ws.Range("H77").Value = 2780.831
$ws.Range("I77").Value = 2503
$ws.Range("K77").Value = 12515
$ws.Range("M77").Value = -8147
$ws.Range("H116").Value = 1139.5
$ws.Range("I116").Value = 1082.7028
$ws.Range("K116").Value = 1082.7028
$ws.Range("M116").Value = 1211.2972
$ws.Range("H132").Value = 3847.3333
$ws.Range("I132").Value = 3224.76
$ws.Range("K132").Value = 9674.280000000001
$ws.Range("M132").Value = -7144.280000000001
$ws.Range("H136").Value = 21657.691
$ws.Range("I136").Value = 25222.143
$ws.Range("K136").Value = 75666.429
$ws.Range("M136").Value = -73116.429

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1139.5
$ws.Range("I3").Value = 1082.7028
$ws.Range("K3").Value = 1082.7028
$ws.Range("M3").Value = -968.7028
$ws.Range("H80").Value = 4057.3
$ws.Range("I80").Value = 6693.75
$ws.Range("K80").Value = 6693.75
$ws.Range("M80").Value = -5695.75
$ws.Range("H83").Value = 4057.3
$ws.Range("I83").Value = 6693.75
$ws.Range("K83").Value = 33468.75
$ws.Range("M83").Value = -28476.75
$ws.Range("H107").Value = 1495.6578
$ws.Range("I107").Value = 1495.6578
$ws.Range("K107").Value = 1495.6578
$ws.Range("M107").Value = 424.3422
$ws.Range("H126").Value = 83225
$ws.Range("J126").Value = 83225
$ws.Range("L126").Value = 83225
$ws.Range("N126").Value = -93105
$ws.Range("H134").Value = 3246.8438
$ws.Range("I134").Value = 2717.7256
$ws.Range("K134").Value = 8153.176800000001
$ws.Range("M134").Value = -5618.176800000001
$ws.Range("H135").Value = 76664.336
$ws.Range("J135").Value = 76664.336
$ws.Range("L135").Value = 76664.336
$ws.Range("N135").Value = -86804.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 303.25
$ws.Range("I22").Value = 281.2
$ws.Range("K22").Value = 281.2
$ws.Range("M22").Value = 68.80000000000001
$ws.Range("H31").Value = 5508.365
$ws.Range("I31").Value = 2077
$ws.Range("J31").Value = 8939.73
$ws.Range("K31").Value = 2077
$ws.Range("L31").Value = 8939.73
$ws.Range("M31").Value = -1782
$ws.Range("N31").Value = -9529.73
$ws.Range("H34").Value = 5508.365
$ws.Range("I34").Value = 2077
$ws.Range("J34").Value = 8939.73
$ws.Range("K34").Value = 2077
$ws.Range("L34").Value = 8939.73
$ws.Range("M34").Value = -1875
$ws.Range("N34").Value = -9343.73
$ws.Range("H58").Value = 3616.6667
$ws.Range("I58").Value = 2200
$ws.Range("K58").Value = 2200
$ws.Range("M58").Value = -1997
$ws.Range("H132").Value = 22729914
$ws.Range("I132").Value = 2643.25
$ws.Range("J132").Value = 83335970
$ws.Range("K132").Value = 7929.75
$ws.Range("L132").Value = 250007910
$ws.Range("M132").Value = -5399.75
$ws.Range("N132").Value = -250012970
$ws.Range("H134").Value = 2289
$ws.Range("I134").Value = 2203.5789
$ws.Range("K134").Value = 6610.736699999999
$ws.Range("M134").Value = -4075.736699999999
$ws.Range("H136").Value = 3616.6667
$ws.Range("I136").Value = 2200
$ws.Range("K136").Value = 6600
$ws.Range("M136").Value = -4050

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 710.1
$ws.Range("I98").Value = 600
$ws.Range("J98").Value = 722.3333
$ws.Range("K98").Value = 1800
$ws.Range("L98").Value = 2166.9999
$ws.Range("M98").Value = -302
$ws.Range("N98").Value = -5162.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 3125.75
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 3125.75
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 3125.75
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -3357.75
$ws.Range("H97").Value = 1850.4375
$ws.Range("I97").Value = 281.95
$ws.Range("J97").Value = 4464.5835
$ws.Range("K97").Value = 281.95
$ws.Range("L97").Value = 4464.5835
$ws.Range("M97").Value = 214.05
$ws.Range("N97").Value = -5456.5835
$ws.Range("H132").Value = 391130.34
$ws.Range("I132").Value = 804757.75
$ws.Range("J132").Value = 3354.6875
$ws.Range("K132").Value = 2414273.25
$ws.Range("L132").Value = 10064.0625
$ws.Range("M132").Value = -2411743.25
$ws.Range("N132").Value = -15124.0625
$ws.Range("H133").Value = 88649.664
$ws.Range("J133").Value = 88649.664
$ws.Range("L133").Value = 88649.664
$ws.Range("N133").Value = -98769.664

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 5003
$ws.Range("I12").Value = 5003
$ws.Range("K12").Value = 5003
$ws.Range("M12").Value = -4833
$ws.Range("H40").Value = 2920.4167
$ws.Range("I40").Value = 2000
$ws.Range("K40").Value = 2000
$ws.Range("M40").Value = -1864
$ws.Range("H122").Value = 52102.617
$ws.Range("J122").Value = 115785.89
$ws.Range("L122").Value = 347357.67
$ws.Range("N122").Value = -352257.67
$ws.Range("H132").Value = 788532
$ws.Range("I132").Value = 1050003.2
$ws.Range("K132").Value = 3150009.6
$ws.Range("M132").Value = -3147479.6
$ws.Range("H136").Value = 4982.3
$ws.Range("I136").Value = 4273.7085
$ws.Range("J136").Value = 7816.6665
$ws.Range("K136").Value = 12821.1255
$ws.Range("L136").Value = 23449.9995
$ws.Range("M136").Value = -10271.1255
$ws.Range("N136").Value = -28549.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3491.64
$ws.Range("I107").Value = 947.53845
$ws.Range("K107").Value = 2842.61535
$ws.Range("M107").Value = -922.61535
$ws.Range("H112").Value = 51443.25
$ws.Range("J112").Value = 51443.25
$ws.Range("L112").Value = 51443.25
$ws.Range("N112").Value = -54397.25
$ws.Range("H132").Value = 518505.06
$ws.Range("I132").Value = 770331.2
$ws.Range("K132").Value = 2310993.6
$ws.Range("M132").Value = -2308463.6
